$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (Strike#) values for column G, rows 2-40, replacing the old values
$newK = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 4
    9  = 2
    10 = 4
    11 = 1
    12 = 5
    13 = 4
    14 = 2
    15 = 3
    16 = 4
    17 = 4
    18 = 4
    19 = 2
    20 = 6
    21 = 3
    22 = 8
    23 = 2
    24 = 7
    25 = 3
    26 = 8
    27 = 4
    28 = 5
    29 = 4
    30 = 3
    31 = 1
    32 = 8
    33 = 4
    34 = 5
    35 = 3
    36 = 3
    37 = 5
    38 = 3
    39 = 3
    40 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
